# LDLC price tracker: append a new snapshot column.
# A new timestamped column is inserted right before the existing "nom" column
# (which, together with "url_produit", always stays as the last two columns).
# For rows that already carry a price series (rows 2-80), the newly inserted
# cell simply repeats the most recent known price (the value currently in the
# last timestamp column, DK). Rows that have no price data yet (81-206) get
# an empty cell in the new column, just like all the other timestamp columns
# on those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DL ("nom" shifts to DM, "url_produit" shifts to DN).
$ws.Columns("DL:DL").Insert()

# New header timestamp, continuing the existing sequence in row 1.
$ws.Range("DL1").Value = "2026-02-01 23:15:25"

# Carry forward the latest known price (column DK) into the new column for
# every product row that already has pricing history.
for ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, "DL").Value = $ws.Cells.Item($r, "DK").Value()
}
